$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Duplicate the existing "2022-Q3" sheet (placing the copy right before
#    it) and rename the copy to "2022-Q4". Duplicating - rather than adding
#    a blank sheet - carries over the sheet's page setup / outline props /
#    header-row style automatically, so the new quarter sheet looks exactly
#    like its siblings.
#    Final sheet order: 总计, 2022-Q4, 2022-Q3, 2022-Q2
# ---------------------------------------------------------------------------
$oldQ3 = $wb.Worksheets.Item("2022-Q3")
$oldQ3.Copy($oldQ3)

$q4 = $wb.Worksheets.Item("2022-Q3 (2)")
$q4.Name = "2022-Q4"

# ---------------------------------------------------------------------------
# 2. Overwrite "2022-Q4" with the new fund-holding data (header row is
#    already correct, copied from "2022-Q3"). Columns B/D/E/F/G hold
#    numeric-looking values that must stay TEXT (matching the source data),
#    so force the cell format to Text before assigning, then drop back to
#    the default "Normal" style (keeps the value text without leaving a
#    stray custom number format behind).
# ---------------------------------------------------------------------------
$q4 = $wb.Worksheets.Item("2022-Q4")

$q4.Range("A2").Value = 0
$q4.Range("B2").NumberFormat = "@"
$q4.Range("B2").Value = "011351"
$q4.Range("B2").Style = "Normal"
$q4.Range("C2").Value = "金鹰年年邮益一年持有期混合A"
$q4.Range("D2").NumberFormat = "@"
$q4.Range("D2").Value = "3.04"
$q4.Range("D2").Style = "Normal"
$q4.Range("E2").NumberFormat = "@"
$q4.Range("E2").Value = "39.17"
$q4.Range("E2").Style = "Normal"
$q4.Range("F2").NumberFormat = "@"
$q4.Range("F2").Value = "0.75"
$q4.Range("F2").Style = "Normal"
$q4.Range("G2").NumberFormat = "@"
$q4.Range("G2").Value = "0.0228"
$q4.Range("G2").Style = "Normal"
$q4.Range("H2").Value = 8

$q4.Range("A3").Value = 1
$q4.Range("B3").NumberFormat = "@"
$q4.Range("B3").Value = "011352"
$q4.Range("B3").Style = "Normal"
$q4.Range("C3").Value = "金鹰年年邮益一年持有期混合C"
$q4.Range("D3").NumberFormat = "@"
$q4.Range("D3").Value = "0.23"
$q4.Range("D3").Style = "Normal"
$q4.Range("E3").NumberFormat = "@"
$q4.Range("E3").Value = "39.17"
$q4.Range("E3").Style = "Normal"
$q4.Range("F3").NumberFormat = "@"
$q4.Range("F3").Value = "0.75"
$q4.Range("F3").Style = "Normal"
$q4.Range("G3").NumberFormat = "@"
$q4.Range("G3").Value = "0.0017"
$q4.Range("G3").Style = "Normal"
$q4.Range("H3").Value = 8

# ---------------------------------------------------------------------------
# 3. Update the "总计" (summary) sheet: shift the existing two data rows
#    down by one and insert the new quarter's totals at the top.
# ---------------------------------------------------------------------------
$tot = $wb.Worksheets.Item("总计")

# New row 4 (2022-Q2 / 0.18, previously row 3) - copy A3's style first so
# the new A4 cell picks up the same bold/bordered formatting.
$tot.Range("A3").Copy($tot.Range("A4"))
$tot.Range("A4").Value = 2
$tot.Range("B4").Value = "2022-Q2"
$tot.Range("C4").Value = 2
$tot.Range("D4").Value = 0.18

# Row 3 becomes 2022-Q3 / 0.01 (previously 2022-Q2 / 0.18).
$tot.Range("B3").Value = "2022-Q3"
$tot.Range("D3").Value = 0.01

# Row 2 becomes 2022-Q4 / 0.02 (previously 2022-Q3 / 0.01).
$tot.Range("B2").Value = "2022-Q4"
$tot.Range("D2").Value = 0.02

# ---------------------------------------------------------------------------
# 4. Restore the originally-active tab ("2022-Q2" was the selected sheet in
#    the source workbook).
# ---------------------------------------------------------------------------
$wb.Worksheets.Item("2022-Q2").Select()
